{"js": "// Update the date line and the 20x5 arithmetic-problem table in place.\n// New values are applied positionally (row-major) to exactly match the\n// target diff, since a few problems (e.g. \"54+44=\") occur more than once\n// with different replacements depending on position.\n\nconst body = context.document.body;\n\n// --- 1) Update the date paragraph (\"2025-03-29 Saturday\" -> \"2025-03-30 Sunday\") ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The date line is the first paragraph in the document (above the table).\n// insertText(..., replace) on the paragraph's own range swaps the text\n// content while keeping the existing run formatting (Arial, sz 30).\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2025-03-30 Sunday\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2) Update every cell of the table, in row-major order ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst newValues = [\n  [\"77-76=\", \"70-39=\", \"41+46=\", \"62-18=\", \"90-15=\"],\n  [\"25-5=\", \"52+15=\", \"29+18=\", \"71+26=\", \"10+57=\"],\n  [\"26+20=\", \"91-69=\", \"90-25=\", \"48+3=\", \"57-16=\"],\n  [\"95-85=\", \"64-54=\", \"13+55=\", \"71-18=\", \"84-13=\"],\n  [\"63+35=\", \"55+39=\", \"44+36=\", \"3+2=\", \"20+29=\"],\n  [\"57-0=\", \"30-13=\", \"19+26=\", \"66-6=\", \"97-62=\"],\n  [\"72-1=\", \"38-25=\", \"25+32=\", \"61+20=\", \"87-64=\"],\n  [\"88-75=\", \"93-31=\", \"22+19=\", \"8+60=\", \"58+30=\"],\n  [\"35+57=\", \"67-14=\", \"35-4=\", \"26-21=\", \"84-62=\"],\n  [\"63-33=\", \"30+9=\", \"18+6=\", \"73-36=\", \"26+0=\"],\n  [\"0+47=\", \"4+10=\", \"86-75=\", \"47-24=\", \"28+21=\"],\n  [\"81-4=\", \"43-22=\", \"24+57=\", \"51-2=\", \"80-73=\"],\n  [\"64-46=\", \"12+24=\", \"73-9=\", \"87-41=\", \"79-27=\"],\n  [\"81-80=\", \"22+12=\", \"97-52=\", \"78-6=\", \"32+37=\"],\n  [\"20+66=\", \"26+32=\", \"13+81=\", \"16+23=\", \"38+14=\"],\n  [\"82-72=\", \"94-71=\", \"23+65=\", \"81-36=\", \"13+37=\"],\n  [\"84+11=\", \"5+42=\", \"58+4=\", \"53-37=\", \"24+39=\"],\n  [\"13+41=\", \"34+7=\", \"15+14=\", \"0+6=\", \"34+36=\"],\n  [\"81-40=\", \"92-3=\", \"76-27=\", \"43-20=\", \"69-52=\"],\n  [\"6+62=\", \"17+57=\", \"53-26=\", \"96-32=\", \"15+33=\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the date line and the 20x5 arithmetic-problem table in place.\n# New values are applied positionally (row by row, left to right) to\n# exactly match the target diff, since a few problems (e.g. \"54+44=\")\n# occur more than once with different replacements depending on position.\n\n$d = $word.ActiveDocument\n\n# --- 1) Update the date paragraph (\"2025-03-29 Saturday\" -> \"2025-03-30 Sunday\") ---\n$dateParagraph = $d.Paragraphs.Item(1)\n$dateParagraph.Range.Text = \"2025-03-30 Sunday\"\n\n# --- 2) Update every cell of the table, in row-major order ---\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"77-76=\", \"70-39=\", \"41+46=\", \"62-18=\", \"90-15=\"),\n    @(\"25-5=\", \"52+15=\", \"29+18=\", \"71+26=\", \"10+57=\"),\n    @(\"26+20=\", \"91-69=\", \"90-25=\", \"48+3=\", \"57-16=\"),\n    @(\"95-85=\", \"64-54=\", \"13+55=\", \"71-18=\", \"84-13=\"),\n    @(\"63+35=\", \"55+39=\", \"44+36=\", \"3+2=\", \"20+29=\"),\n    @(\"57-0=\", \"30-13=\", \"19+26=\", \"66-6=\", \"97-62=\"),\n    @(\"72-1=\", \"38-25=\", \"25+32=\", \"61+20=\", \"87-64=\"),\n    @(\"88-75=\", \"93-31=\", \"22+19=\", \"8+60=\", \"58+30=\"),\n    @(\"35+57=\", \"67-14=\", \"35-4=\", \"26-21=\", \"84-62=\"),\n    @(\"63-33=\", \"30+9=\", \"18+6=\", \"73-36=\", \"26+0=\"),\n    @(\"0+47=\", \"4+10=\", \"86-75=\", \"47-24=\", \"28+21=\"),\n    @(\"81-4=\", \"43-22=\", \"24+57=\", \"51-2=\", \"80-73=\"),\n    @(\"64-46=\", \"12+24=\", \"73-9=\", \"87-41=\", \"79-27=\"),\n    @(\"81-80=\", \"22+12=\", \"97-52=\", \"78-6=\", \"32+37=\"),\n    @(\"20+66=\", \"26+32=\", \"13+81=\", \"16+23=\", \"38+14=\"),\n    @(\"82-72=\", \"94-71=\", \"23+65=\", \"81-36=\", \"13+37=\"),\n    @(\"84+11=\", \"5+42=\", \"58+4=\", \"53-37=\", \"24+39=\"),\n    @(\"13+41=\", \"34+7=\", \"15+14=\", \"0+6=\", \"34+36=\"),\n    @(\"81-40=\", \"92-3=\", \"76-27=\", \"43-20=\", \"69-52=\"),\n    @(\"6+62=\", \"17+57=\", \"53-26=\", \"96-32=\", \"15+33=\")\n)\n\n$rowCount = $newValues.Count\nfor ($r = 0; $r -lt $rowCount; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n"}
